{"js": "// Replace each three-digit-division expression in the worksheet table\n// with its updated value, matching the exact text of each cell.\nconst replacements = [\n  [\"275\u00f76=\", \"124\u00f74=\"],\n  [\"446\u00f73=\", \"740\u00f78=\"],\n  [\"720\u00f74=\", \"876\u00f75=\"],\n  [\"532\u00f77=\", \"501\u00f75=\"],\n  [\"836\u00f76=\", \"406\u00f74=\"],\n  [\"154\u00f72=\", \"764\u00f75=\"],\n  [\"353\u00f78=\", \"259\u00f77=\"],\n  [\"941\u00f76=\", \"525\u00f77=\"],\n  [\"342\u00f79=\", \"202\u00f76=\"],\n  [\"959\u00f72=\", \"754\u00f73=\"],\n  [\"337\u00f77=\", \"420\u00f76=\"],\n  [\"407\u00f75=\", \"151\u00f73=\"],\n  [\"754\u00f74=\", \"960\u00f72=\"],\n  [\"716\u00f78=\", \"565\u00f77=\"],\n  [\"624\u00f77=\", \"188\u00f79=\"],\n  [\"797\u00f72=\", \"532\u00f79=\"],\n  [\"518\u00f79=\", \"208\u00f74=\"],\n  [\"116\u00f76=\", \"402\u00f78=\"],\n  [\"115\u00f72=\", \"908\u00f72=\"],\n  [\"594\u00f74=\", \"965\u00f72=\"],\n  [\"572\u00f72=\", \"550\u00f73=\"],\n  [\"871\u00f78=\", \"946\u00f77=\"],\n  [\"833\u00f73=\", \"118\u00f79=\"],\n  [\"838\u00f78=\", \"637\u00f78=\"],\n  [\"825\u00f73=\", \"738\u00f77=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update each three-digit-division expression in the worksheet table\n# to its new value by finding the exact old text and replacing it.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"275\u00f76=\", \"124\u00f74=\"),\n    @(\"446\u00f73=\", \"740\u00f78=\"),\n    @(\"720\u00f74=\", \"876\u00f75=\"),\n    @(\"532\u00f77=\", \"501\u00f75=\"),\n    @(\"836\u00f76=\", \"406\u00f74=\"),\n    @(\"154\u00f72=\", \"764\u00f75=\"),\n    @(\"353\u00f78=\", \"259\u00f77=\"),\n    @(\"941\u00f76=\", \"525\u00f77=\"),\n    @(\"342\u00f79=\", \"202\u00f76=\"),\n    @(\"959\u00f72=\", \"754\u00f73=\"),\n    @(\"337\u00f77=\", \"420\u00f76=\"),\n    @(\"407\u00f75=\", \"151\u00f73=\"),\n    @(\"754\u00f74=\", \"960\u00f72=\"),\n    @(\"716\u00f78=\", \"565\u00f77=\"),\n    @(\"624\u00f77=\", \"188\u00f79=\"),\n    @(\"797\u00f72=\", \"532\u00f79=\"),\n    @(\"518\u00f79=\", \"208\u00f74=\"),\n    @(\"116\u00f76=\", \"402\u00f78=\"),\n    @(\"115\u00f72=\", \"908\u00f72=\"),\n    @(\"594\u00f74=\", \"965\u00f72=\"),\n    @(\"572\u00f72=\", \"550\u00f73=\"),\n    @(\"871\u00f78=\", \"946\u00f77=\"),\n    @(\"833\u00f73=\", \"118\u00f79=\"),\n    @(\"838\u00f78=\", \"637\u00f78=\"),\n    @(\"825\u00f73=\", \"738\u00f77=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
